$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program_choosing")

$ws.Range("A4").Value = "TUM Business Informatics"
$ws.Range("B4").Value = "Yes"

$ws.Range("A5").Value = "Tuebingen Machine Learning"
$ws.Range("B5").Value = "Yes"

$ws.Range("A6").Select()
